$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply an AutoFilter on column B ("Grupa") keeping only rows where the
# value equals 244 - this hides every non-matching data row (2-73, 106).
$ws.Range("A1:J106").AutoFilter(2, @("244"))

# Zoom the view to 80% and move the selection to J72, matching the
# author's final on-screen state.
$excel.ActiveWindow.Zoom = 80
$ws.Range("J72").Select()
